$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(245, 44319, 1, 4, 95.30617107457708),
    @(246, 44320, 0, 4, 95.30617107457708),
    @(247, 44321, 0, 4, 95.30617107457708)
)

$srcDateCell = $ws.Cells.Item(244, 1)

foreach ($row in $data) {
    $r = $row[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $srcDateCell.Copy()
    $dateCell.PasteSpecial(-4122)
    $dateCell.Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
